$d = $word.ActiveDocument

$d.Content.Find.Execute("551÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "534÷5=", 2) | Out-Null
$d.Content.Find.Execute("281÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "919÷7=", 2) | Out-Null
$d.Content.Find.Execute("199÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "646÷4=", 2) | Out-Null
$d.Content.Find.Execute("147÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "350÷5=", 2) | Out-Null
$d.Content.Find.Execute("919÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "227÷2=", 2) | Out-Null
$d.Content.Find.Execute("992÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "407÷8=", 2) | Out-Null
$d.Content.Find.Execute("270÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "273÷9=", 2) | Out-Null
$d.Content.Find.Execute("990÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "386÷6=", 2) | Out-Null
$d.Content.Find.Execute("466÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "262÷7=", 2) | Out-Null
$d.Content.Find.Execute("590÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "245÷6=", 2) | Out-Null
$d.Content.Find.Execute("618÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "921÷8=", 2) | Out-Null
$d.Content.Find.Execute("688÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "320÷8=", 2) | Out-Null
$d.Content.Find.Execute("598÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "664÷7=", 2) | Out-Null
$d.Content.Find.Execute("495÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "211÷4=", 2) | Out-Null
$d.Content.Find.Execute("222÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "423÷4=", 2) | Out-Null
$d.Content.Find.Execute("122÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "928÷8=", 2) | Out-Null
$d.Content.Find.Execute("468÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "188÷8=", 2) | Out-Null
$d.Content.Find.Execute("661÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "180÷2=", 2) | Out-Null
$d.Content.Find.Execute("237÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "479÷8=", 2) | Out-Null
$d.Content.Find.Execute("206÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "444÷7=", 2) | Out-Null
$d.Content.Find.Execute("732÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "230÷4=", 2) | Out-Null
$d.Content.Find.Execute("158÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "336÷6=", 2) | Out-Null
$d.Content.Find.Execute("497÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "883÷6=", 2) | Out-Null
$d.Content.Find.Execute("627÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "497÷6=", 2) | Out-Null
$d.Content.Find.Execute("537÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "679÷8=", 2) | Out-Null
